# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted right
# before the existing "In Advance" column (previously column N), so that
# the schedule has an extra spacer/variable-instalment column.  This pushes
# "In Advance" -> O, the blank spacer that used to be O -> P, and
# "Outstanding" -> Q.
#
# Also update which sheet/cell is active (the workbook was left with the
# "Repayment Schedule" tab selected, cell S9, instead of the "Transactions"
# tab).

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (the "In Advance" column),
# shifting it (and everything after) one column to the right.
$wsSchedule.Columns.Item(14).Insert()

# Give the newly inserted column a sensible width (matches the width of
# the neighbouring "In Advance"/Principal columns).
$wsSchedule.Columns.Item(14).ColumnWidth = 10.5

# Make "Repayment Schedule" the active sheet/tab, with cell S9 selected,
# instead of "Transactions".
$wsSchedule.Select()
$wsSchedule.Range("S9").Select()
